$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 0.55
$ws.Range("P2").Value = 0.15
$ws.Range("S2").Value = 0.1

# Row 3
$ws.Range("P3").Value = 0.7272727272727273
$ws.Range("S3").Value = 0.2727272727272727

# Row 4
$ws.Range("S4").Value = 1

# Row 6
$ws.Range("B6").Value = 0.12
$ws.Range("J6").Value = 0.08
$ws.Range("O6").Value = 0.04
$ws.Range("Q6").Value = 0.12
$ws.Range("R6").Value = 0.08
$ws.Range("S6").Value = 0.56

# Row 7
$ws.Range("B7").Value = 0.08695652173913043
$ws.Range("F7").Value = 0.1304347826086956
$ws.Range("Q7").Value = 0.04347826086956522
$ws.Range("R7").Value = 0.08695652173913043
$ws.Range("S7").Value = 0.6521739130434783

# Row 8
$ws.Range("B8").Value = 0.03658536585365853
$ws.Range("F8").Value = 0.07317073170731707
$ws.Range("J8").Value = 0.0975609756097561
$ws.Range("O8").Value = 0.02439024390243903
$ws.Range("Q8").Value = 0.0975609756097561
$ws.Range("R8").Value = 0.07317073170731707
$ws.Range("S8").Value = 0.5975609756097561

# Row 9
$ws.Range("B9").Value = 0.08571428571428572
$ws.Range("D9").Value = 0.02857142857142857
$ws.Range("F9").Value = 0.08571428571428572
$ws.Range("J9").Value = 0.08571428571428572
$ws.Range("O9").Value = 0.05714285714285714
$ws.Range("Q9").Value = 0.1428571428571428
$ws.Range("R9").Value = 0.1428571428571428
$ws.Range("S9").Value = 0.3714285714285714

# Row 10
$ws.Range("B10").Value = 0.04424778761061947
$ws.Range("F10").Value = 0.07079646017699115
$ws.Range("J10").Value = 0.1150442477876106
$ws.Range("O10").Value = 0.008849557522123894
$ws.Range("Q10").Value = 0.1769911504424779
$ws.Range("R10").Value = 0.06194690265486726
$ws.Range("S10").Value = 0.5221238938053098

# Row 11
$ws.Range("G11").Value = 0.125
$ws.Range("J11").Value = 0.03125
$ws.Range("K11").Value = 0.15625
$ws.Range("L11").Value = 0.625
$ws.Range("S11").Value = 0.0625

# Row 12
$ws.Range("G12").Value = 0.85
$ws.Range("J12").Value = 0.15

# Row 13
$ws.Range("G13").Value = 0.5
$ws.Range("J13").Value = 0.25
$ws.Range("S13").Value = 0.25

# Row 15
$ws.Range("H15").Value = 0.09523809523809523
$ws.Range("I15").Value = 0.1428571428571428
$ws.Range("J15").Value = 0.2380952380952381
$ws.Range("S15").Value = 0.5238095238095238

# Row 16
$ws.Range("H16").Value = 0.3636363636363636
$ws.Range("I16").Value = 0.09090909090909091
$ws.Range("J16").Value = 0.1818181818181818
$ws.Range("K16").Value = 0.1818181818181818
$ws.Range("O16").Value = 0.09090909090909091
$ws.Range("S16").Value = 0.09090909090909091

# Row 17
$ws.Range("H17").Value = 0.2432432432432433
$ws.Range("I17").Value = 0.2162162162162162
$ws.Range("J17").Value = 0.3783783783783784
$ws.Range("S17").Value = 0.1621621621621622

# Row 18
$ws.Range("H18").Value = 0.3809523809523809
$ws.Range("I18").Value = 0.1428571428571428
$ws.Range("J18").Value = 0.1428571428571428
$ws.Range("K18").Value = 0.09523809523809523
$ws.Range("O18").Value = 0.04761904761904762
$ws.Range("S18").Value = 0.1904761904761905

# Row 19
$ws.Range("F19").Value = 0.009478672985781991
$ws.Range("H19").Value = 0.2796208530805687
$ws.Range("I19").Value = 0.0995260663507109
$ws.Range("J19").Value = 0.2843601895734597
$ws.Range("K19").Value = 0.1090047393364929
$ws.Range("M19").Value = 0.02369668246445497
$ws.Range("O19").Value = 0.04265402843601896
$ws.Range("S19").Value = 0.1516587677725119
